$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Dennis Schröder"
$ws.Range("B2").Value = "PG,SG"
$ws.Range("C2").Value = "Golden State Warriors"
$ws.Range("A3").Value = "Amen Thompson"
$ws.Range("B3").Value = "SG,SF"
$ws.Range("C3").Value = "Houston Rockets"
$ws.Range("A4").Value = "Fred VanVleet"
$ws.Range("B4").Value = "PG"
$ws.Range("C4").Value = "Houston Rockets"
$ws.Range("A5").Value = "Dillon Brooks"
$ws.Range("B5").Value = "SG,SF"
$ws.Range("C5").Value = "Houston Rockets"
$ws.Range("A6").Value = "Anthony Edwards"
$ws.Range("B6").Value = "SG,SF"
$ws.Range("C6").Value = "Minnesota Timberwolves"
$ws.Range("A7").Value = "Jayson Tatum"
$ws.Range("B7").Value = "SF,PF"
$ws.Range("C7").Value = "Boston Celtics"
$ws.Range("A8").Value = "Grayson Allen"
$ws.Range("B8").Value = "PG,SG,SF"
$ws.Range("C8").Value = "Phoenix Suns"
$ws.Range("A9").Value = "Zion Williamson"
$ws.Range("B9").Value = "PF,C"
$ws.Range("C9").Value = "New Orleans Pelicans"
$ws.Range("A10").Value = "Jaren Jackson Jr."
$ws.Range("B10").Value = "PF,C"
$ws.Range("C10").Value = "Memphis Grizzlies"
$ws.Range("A11").Value = "Ivica Zubac"
$ws.Range("B11").Value = "C"
$ws.Range("C11").Value = "LA Clippers"
$ws.Range("A12").Value = "Andrew Wiggins"
$ws.Range("B12").Value = "SF,PF"
$ws.Range("C12").Value = "Golden State Warriors"
$ws.Range("A13").Value = "Scoot Henderson"
$ws.Range("B13").Value = "PG"
$ws.Range("C13").Value = "Portland Trail Blazers"
$ws.Range("A14").Value = "Anfernee Simons"
$ws.Range("B14").Value = "PG,SG"
$ws.Range("C14").Value = "Portland Trail Blazers"
$ws.Range("A15").Value = "James Harden"
$ws.Range("B15").Value = "PG,SG"
$ws.Range("C15").Value = "LA Clippers"
$ws.Range("A16").Value = "Paul George"
$ws.Range("B16").Value = "SG,SF,PF"
$ws.Range("C16").Value = "Philadelphia 76ers"
$ws.Range("A17").Value = "Jonathan Kuminga"
$ws.Range("B17").Value = "SF,PF"
$ws.Range("C17").Value = "Golden State Warriors"
$ws.Range("A18").Value = "Jerami Grant"
$ws.Range("B18").Value = "SG,SF,PF"
$ws.Range("C18").Value = "Portland Trail Blazers"
$ws.Range("A19").Value = "Giannis Antetokounmpo"
$ws.Range("B19").Value = "PF,C"
$ws.Range("C19").Value = "Milwaukee Bucks"
